$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Rename header cells in row 1 (A1:U1) to the new "<formatversion>" suffix scheme.
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Step 2: Freeze the header row (split at row 1, so row 2 is the first scrollable row).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Step 3: Convert the range A1:U93 into an Excel Table (ListObject) named "Table1".
$range = $ws.Range("A1:U93")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
